# Refresh the cryptocurrency price/volume snapshot (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.649.05'
$ws.Range('E2').Value = '  -4.78%  '

$ws.Range('D3').Value = '1.946.04'
$ws.Range('E3').Value = '  -4.74%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = '''241.18'
$ws.Range('E5').Value = '  -4.49%  '

$ws.Range('D6').Value = '''0.619'
$ws.Range('E6').Value = '  -4.80%  '

$ws.Range('D7').Value = '''60.38'
$ws.Range('E7').Value = '  -7.60%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').Value = '''0.363'
$ws.Range('E9').Value = '  -3.26%  '

$ws.Range('D10').Value = '''55.86'
$ws.Range('E10').Value = '  -5.51%  '

$ws.Range('D11').Value = '''0.0782'
$ws.Range('E11').Value = '  +3.51%  '

$ws.Range('E12').Value = '  -1.74%  '

$ws.Range('D13').Value = '''0.845'
$ws.Range('E13').Value = '  -6.90%  '

$ws.Range('D14').Value = '''13.75'
$ws.Range('E14').Value = '  -9.20%  '

$ws.Range('D15').Value = '2.242.23'
$ws.Range('E15').Value = '  -4.11%  '

$ws.Range('D16').Value = '''21.17'
$ws.Range('E16').Value = '  +1.04%  '

$ws.Range('E17').Value = '  -4.59%  '

$ws.Range('D18').Value = '1.960.65'
$ws.Range('E18').Value = '  -3.86%  '

$ws.Range('D19').Value = '35.516.47'
$ws.Range('E19').Value = '  -4.99%  '

$ws.Range('D20').Value = '''70.13'
$ws.Range('E20').Value = '  -4.16%  '

$ws.Range('E21').Value = '  -4.07%  '

$ws.Range('D22').Value = '''236.33'
$ws.Range('E22').Value = '  -0.09%  '

$ws.Range('D23').Value = '''5.12'
$ws.Range('E23').Value = '  -4.39%  '

$ws.Range('E24').Value = '  -0.05%  '

$ws.Range('D25').Value = '''2.49'
$ws.Range('E25').Value = '  -9.44%  '

$ws.Range('D26').Value = '''2.27'
$ws.Range('E26').Value = '  -3.20%  '

$ws.Range('D27').Value = '''9.57'
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').Value = '''157.91'
$ws.Range('E28').Value = '  -4.86%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''19.50'
$ws.Range('E29').Value = '  -1.89%  '

$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '''0.128'
$ws.Range('E30').Value = '  +14.28%  '

$ws.Range('E31').Value = '  -2.91%  '

$ws.Range('D32').Value = '''4.79'
$ws.Range('E32').Value = '  -8.19%  '

$ws.Range('E33').Value = '  -7.98%  '

$ws.Range('D34').Value = '''0.0607'
$ws.Range('E34').Value = '  -1.50%  '

$ws.Range('D35').Value = '''4.31'
$ws.Range('E35').Value = '  -9.04%  '

$ws.Range('D36').Value = '''6.17'
$ws.Range('E36').Value = '  +3.45%  '

$ws.Range('E37').Value = '  +0.05%  '

$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''2.27'
$ws.Range('E38').Value = '  -7.75%  '

$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').Value = '''1.83'
$ws.Range('E39').Value = '  +1.34%  '

$ws.Range('D40').Value = '''3.05'
$ws.Range('E40').Value = '  +11.41%  '

$ws.Range('D41').Value = '''0.0971'
$ws.Range('E41').Value = '  -6.91%  '

$ws.Range('E42').Value = '  -2.42%  '

$ws.Range('D43').Value = '''2.79'
$ws.Range('E43').Value = '  -4.81%  '

$ws.Range('D44').Value = '''0.0209'
$ws.Range('E44').Value = '  -4.60%  '

$ws.Range('D45').Value = '''1.07'
$ws.Range('E45').Value = '  -6.01%  '

$ws.Range('D46').Value = '''91.57'
$ws.Range('E46').Value = '  -4.26%  '

$ws.Range('D47').Value = '''15.79'
$ws.Range('E47').Value = '  -7.91%  '

$ws.Range('D48').Value = '''7.43'
$ws.Range('E48').Value = '  -8.45%  '

$ws.Range('D49').Value = '1.327.09'
$ws.Range('E49').Value = '  -6.50%  '

$ws.Range('E50').Value = '  -7.26%  '

$ws.Range('D51').Value = '2.132.64'
$ws.Range('E51').Value = '  -4.10%  '
